$wb = $excel.ActiveWorkbook

$wsTopCities = $wb.Worksheets.Item("TopCities")
$wsInput     = $wb.Worksheets.Item("input")
$wsHospital  = $wb.Worksheets.Item("hospitalName")

# ---------------------------------------------------------------------
# "input" sheet: replace the sample rows of data (rohan/siva/krihsna)
# with the new people (Akash/Hada/Sharma). The order in which new
# string values are assigned below matches the order they need to be
# introduced into the workbook's shared-string table.
# ---------------------------------------------------------------------
$wsInput.Range("A1").Value = "Akash"
$wsInput.Range("A2").Value = "Hada"
$wsInput.Range("A3").Value = "Sharma"

$wsInput.Range("D1").Value = "akashsingh.4111@gmail.com"
$wsInput.Range("D2").Value = "hada@gmail.com"
$wsInput.Range("D3").Value = "sharma"

$wsInput.Range("C2").Value = "s"

$wsInput.Range("C1").Value = 8317096478

# ---------------------------------------------------------------------
# "hospitalName" sheet: insert a new hospital at the top of the list,
# pushing the existing rows down by one.
# ---------------------------------------------------------------------
$wsHospital.Rows.Item(1).Insert() | Out-Null
$wsHospital.Range("B1").Value = "Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)"

# ---------------------------------------------------------------------
# View / selection state changes: "input" becomes the active sheet,
# and each sheet keeps a different selected cell.
# ---------------------------------------------------------------------
$wsTopCities.Activate() | Out-Null
$wsTopCities.Range("D9").Select() | Out-Null

$wsInput.Activate() | Out-Null
$wsInput.Range("C2").Select() | Out-Null

$wsHospital.Activate() | Out-Null
$wsHospital.Range("F9").Select() | Out-Null

$wsInput.Activate() | Out-Null
